# Update the "seas_id" column (column O) values to reflect the refreshed
# CSV data pulled from kaggle.com. Each row below maps to its new value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 4444
    3  = 2404
    4  = 3465
    5  = 4444
    6  = 1951
    7  = 4444
    8  = 4444
    9  = 1525
    10 = 2986
    11 = 2986
    12 = 2660
    13 = 3445
    14 = 5041
    15 = 1556
    16 = 1556
    17 = 1556
    18 = 3097
    19 = 3453
    22 = 3625
    25 = 1424
    28 = 1424
    29 = 1424
    30 = 1849
    31 = 1424
    32 = 2039
    33 = 2039
    34 = 1241
    35 = 1312
    36 = 1312
    37 = 3489
    38 = 3649
    39 = 3638
    41 = 3489
    42 = 4612
    43 = 4612
}

foreach ($row in $updates.Keys) {
    $ws.Range("O$row").Value = $updates[$row]
}
